$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "30.588.69"
$ws.Range("E2").Value = "  +1.59%  "

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.47"
$ws.Range("E3").Value = "  +1.51%  "

# Row 4
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "247.27"
$ws.Range("E5").Value = "  +6.02%  "

# Row 6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4761"
$ws.Range("E7").Value = "  +1.52%  "

# Row 8
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2915"
$ws.Range("E8").Value = "  +3.28%  "

# Row 9
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06538"
$ws.Range("E9").Value = "  +1.48%  "

# Row 10
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "22.12"
$ws.Range("E10").Value = "  +5.75%  "

# Row 11
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07737"
$ws.Range("E11").Value = "  +0.04%  "

# Row 12
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "97.61"
$ws.Range("E12").Value = "  +4.66%  "

# Row 13
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7435"
$ws.Range("E13").Value = "  +9.95%  "

# Row 14
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "1.872.74"
$ws.Range("E14").Value = "  +0.81%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.172"

# Row 16
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "275.65"
$ws.Range("E16").Value = "  +3.72%  "

# Row 17
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "30.566.73"
$ws.Range("E17").Value = "  +1.60%  "

# Row 18
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "13.61"
$ws.Range("E18").Value = "  +2.55%  "

# Row 19
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007592"
$ws.Range("E19").Value = "  +0.54%  "

# Row 20
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "1.0000"
$ws.Range("E20").Value = "  -0.09%  "

# Row 21
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "2.124.37"
$ws.Range("E21").Value = "  +0.56%  "

# Row 22
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "5.286"
$ws.Range("E22").Value = "  +3.17%  "

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9997"
$ws.Range("E23").Value = "  -0.16%  "

# Row 24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "6.220"
$ws.Range("E24").Value = "  +2.24%  "

# Row 25
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "9.334"
$ws.Range("E25").Value = "  +0.73%  "

# Row 26
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "163.36"
$ws.Range("E26").Value = "  -1.14%  "

# Row 27
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "18.99"
$ws.Range("E27").Value = "  +2.88%  "

# Row 28
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "1.949"
$ws.Range("E28").Value = "  +3.96%  "

# Row 29
$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = "1.376"
$ws.Range("E29").Value = "  +0.92%  "

# Row 30
$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09991"
$ws.Range("E30").Value = "  +1.86%  "

# Row 31
$ws.Range("D31:E31").NumberFormat = "@"
$ws.Range("D31").Value = "1.521"
$ws.Range("E31").Value = "  +4.92%  "

# Row 32
$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D32").Value = "4.334"
$ws.Range("E32").Value = "  +3.64%  "

# Row 33
$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D33").Value = "4.089"
$ws.Range("E33").Value = "  +3.32%  "

# Row 34
$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04809"
$ws.Range("E34").Value = "  +3.94%  "

# Row 35
$ws.Range("D35:E35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  +1.90%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7042"

# Row 37
$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  +0.17%  "

# Row 38
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01876"
$ws.Range("E38").Value = "  +2.93%  "

# Row 39
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "2.731"
$ws.Range("E39").Value = "  +0.83%  "

# Row 40
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "6.343"
$ws.Range("E40").Value = "  +1.79%  "

# Row 43
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4232"
$ws.Range("E43").Value = "  +5.14%  "

# Row 44
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8412"
$ws.Range("E44").Value = "  +1.30%  "

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9993"
$ws.Range("E45").Value = "  -0.07%  "

# Row 46
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "102.83"
$ws.Range("E46").Value = "  +1.17%  "

# Row 47
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "9.286"
$ws.Range("E47").Value = "  +1.81%  "

# Row 48
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "7.115"
$ws.Range("E48").Value = "  +3.05%  "

# Row 49
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "35.68"
$ws.Range("E49").Value = "  +5.15%  "

# Row 50
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "921.99"
$ws.Range("E50").Value = "  +0.50%  "

# Row 51
$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3911"
$ws.Range("E51").Value = "  +5.20%  "

# Row 41 and 42 swap: RenderToken/Aave -> Aave/RenderToken
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "71.58"
$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "1.965"
$ws.Range("E42").Value = "  +5.78%  "
